$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# Row 2: int_perc_ipt_age0to5 -> scenario_1 (L2) and scenario_2 (M2)
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 100

# Row 3: int_perc_ipt_age5to15 -> scenario_2 (M3)
$ws.Range("M3").Value = 100

# Row 10 (new): int_perc_xpert
$ws.Range("A10").Value = "int_perc_xpert"
$ws.Range("B10").Value = "no"
$ws.Range("N10").Value = 100

# Update the view selection to reflect the active cell after editing
$ws.Range("A15").Select() | Out-Null
